$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number format on Price column cells so values like "1.00" or "596.49"
# are not auto-converted to numeric values by Excel, matching the original inline string type.
$priceCells = @("D2", "D3", "D5", "D6", "D13", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D30", "D31", "D37", "D38", "D47", "D49", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.105.63'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '2.642.26'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '596.49'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').Value = '155.69'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '27.96'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '3.124.28'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '68.145.45'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '2.660.23'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '11.34'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').Value = '362.87'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').Value = '4.41'
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('D22').Value = '4.78'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('D24').Value = '74.85'
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  -4.02%  '
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = '554.41'
$ws.Range('E30').Value = '  -5.06%  '
$ws.Range('D31').Value = '8.01'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').Value = '161.17'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').Value = '19.41'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('E40').Value = '  -3.39%  '
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('E42').Value = '  +4.72%  '
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('D47').Value = '159.14'
$ws.Range('E47').Value = '  +2.11%  '
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('D49').Value = '21.98'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0782'
$ws.Range('E51').Value = '  -0.21%  '

# Clear the temporary text format so cell style matches the original (no explicit style).
foreach ($c in $priceCells) {
    $ws.Range($c).ClearFormats()
}
